$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "57.861.68"
$ws.Range("E2").Formula = "  +2.09%  "

# Row 3
$ws.Range("D3").Formula = "3.060.32"
$ws.Range("E3").Formula = "  +2.52%  "

# Row 5
$ws.Range("D5").Value = "'526.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "  +6.01%  "

# Row 6
$ws.Range("D6").Value = "'143.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "  +5.74%  "

# Row 7
$ws.Range("E7").Formula = "  +0.00%  "

# Row 8
$ws.Range("D8").Value = "'0.448"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Formula = "  +5.39%  "

# Row 9
$ws.Range("D9").Value = "'7.65"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "  +5.88%  "

# Row 10
$ws.Range("E10").Formula = "  +7.20%  "

# Row 11
$ws.Range("E11").Formula = "  +5.53%  "

# Row 12
$ws.Range("E12").Formula = "  +2.07%  "

# Row 13
$ws.Range("D13").Formula = "3.591.83"
$ws.Range("E13").Formula = "  +2.83%  "

# Row 14
$ws.Range("D14").Value = "'27.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "  +7.91%  "

# Row 15
$ws.Range("D15").Value = "'0.0000171"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "  +15.90%  "

# Row 16
$ws.Range("D16").Formula = "57.884.57"
$ws.Range("E16").Formula = "  +2.29%  "

# Row 17
$ws.Range("D17").Value = "'6.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "  +7.24%  "

# Row 18
$ws.Range("D18").Formula = "3.059.77"
$ws.Range("E18").Formula = "  +2.59%  "

# Row 19
$ws.Range("D19").Value = "'13.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "  +6.68%  "

# Row 20
$ws.Range("E20").Formula = "  +5.24%  "

# Row 21
$ws.Range("D21").Value = "'341.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "  +4.76%  "

# Row 23
$ws.Range("E23").Formula = "  +7.03%  "

# Row 24
$ws.Range("E24").Formula = "  +5.22%  "

# Row 25
$ws.Range("B25").Formula = "Kaspa"
$ws.Range("C25").Formula = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "'0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "  +5.08%  "

# Row 26
$ws.Range("B26").Formula = "PEPE"
$ws.Range("C26").Formula = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Formula = "0.0₃0976"
$ws.Range("E26").Formula = "  +8.82%  "

# Row 27
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "  -0.07%  "

# Row 28
$ws.Range("D28").Value = "'6.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Formula = "  +7.51%  "

# Row 29
$ws.Range("D29").Value = "'7.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Formula = "  +9.60%  "

# Row 30
$ws.Range("E30").Formula = "  +6.63%  "

# Row 31
$ws.Range("E31").Formula = "  +7.09%  "

# Row 32
$ws.Range("D32").Value = "'21.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "  +3.72%  "

# Row 33
$ws.Range("D33").Value = "'4.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "  +7.26%  "

# Row 34
$ws.Range("D34").Value = "'157.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Formula = "  +2.94%  "

# Row 35
$ws.Range("D35").Value = "'5.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "  +6.85%  "

# Row 36
$ws.Range("E36").Formula = "  +3.96%  "

# Row 37
$ws.Range("D37").Value = "'26.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "  +12.45%  "

# Row 38
$ws.Range("D38").Value = "'0.0705"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Formula = "  +4.90%  "

# Row 39
$ws.Range("D39").Formula = "3.098.72"
$ws.Range("E39").Formula = "  +2.74%  "

# Row 40
$ws.Range("D40").Value = "'37.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "  +3.45%  "

# Row 41
$ws.Range("D41").Value = "'3.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "  +10.17%  "

# Row 42
$ws.Range("B42").Formula = "FirstDigitalUSD"
$ws.Range("C42").Formula = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "  +0.23%  "

# Row 43
$ws.Range("B43").Formula = "Mantle"
$ws.Range("C43").Formula = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.667"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "  +4.19%  "

# Row 44
$ws.Range("B44").Formula = "Stacks"
$ws.Range("C44").Formula = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Formula = "  +5.28%  "

# Row 45
$ws.Range("D45").Formula = "2.341.45"
$ws.Range("E45").Formula = "  +5.06%  "

# Row 46
$ws.Range("E46").Formula = "  +3.41%  "

# Row 47
$ws.Range("E47").Formula = "  +2.52%  "

# Row 48
$ws.Range("E48").Formula = "  +5.27%  "

# Row 49
$ws.Range("E49").Formula = "  +3.84%  "

# Row 50
$ws.Range("D50").Value = "'20.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "  +5.99%  "

# Row 51
$ws.Range("D51").Value = "'0.0902"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "  +6.23%  "
